# Recreate the "Таблица с данными" layout on a new sheet "График рассеяния",
# but replace the delta column (D) with a percentage-influence formula and
# a new header label "% влияния".

$wb = $excel.ActiveWorkbook

# --- add the new sheet as the LAST tab ------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "График рассеяния"

# --- header row -------------------------------------------------------------
$ws.Range("A1").Value = "Показатели прибыли"
$ws.Range("B1").Value = 2021
$ws.Range("C1").Value = 2022
$ws.Range("D1").Value = "% влияния"

# --- data rows (same figures as "Таблица с данными") -----------------------
$labels = @(
    "Выручка с логистики",
    "Себестоимость",
    "Валовая прибыль",
    "Комерческие расходы",
    "Управленческие расходы",
    "Прибыль от услуг"
)
$values2021 = @(11, 9, 6, 6, 1, 27)
$values2022 = @(12, 8, 7, 2, 1, 29)

$r = 2
for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Range("A$r").Value = $labels[$i]
    $ws.Range("B$r").Value = $values2021[$i]
    $ws.Range("C$r").Value = $values2022[$i]
    $ws.Range("D$r").Formula = "=C$r/B$r*100"
    $r++
}

# --- formatting: match the "locked" cell style used throughout the workbook
$ws.Range("A1:D7").Locked = $true

# --- sheet view / window ------------------------------------------------
$ws.Range("A1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100

# --- page setup, mirroring the other sheets in the workbook ---------------
$ps = $ws.PageSetup
$ps.PrintHeadings = $false
$ps.PrintGridlines = $false
$ps.PaperSize = 9
$ps.Zoom = 100
$ps.FitToPagesWide = 1
$ps.FitToPagesTall = 1
$ps.Orientation = 1
$ps.LeftMargin = 50.45669291338584
$ps.RightMargin = 50.45669291338584
$ps.TopMargin = 54.14173228346456
$ps.BottomMargin = 54.14173228346456
$ps.HeaderMargin = 21.6
$ps.FooterMargin = 21.6
$ps.CenterHeader = ""

# --- make the new sheet the active tab (activeTab becomes 3) --------------
$ws.Activate()
